# ApplicationTracker.xlsx update:
#   - Add "Code" / "Known Vulnerabilities" / "Found by Agent" columns (G/H/I)
#   - Mark altoro_mutual as having known vulns (F/G -> Y), fill in vuln counts
#   - Add a "Known Vulns" summary table below the main table
#   - Adjust selection / column width to match the authored workbook

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New header cells (row 1) - inherit the bold "header" row style automatically
# ---------------------------------------------------------------------------
$ws.Range("G1").Value = "Code"
$ws.Range("H1").Value = "Known Vulnerabilities"
$ws.Range("I1").Value = "Found by Agent"

# Column H (Known Vulnerabilities count) needs to be wide enough to read
$ws.Columns.Item(8).ColumnWidth = 18.15

# ---------------------------------------------------------------------------
# altoro_mutual (row 5) now has a documented vulnerability -> flip F5 N -> Y
# ---------------------------------------------------------------------------
$ws.Range("F2").Value = "Y"
$ws.Range("F2").Style = "Good"

# ---------------------------------------------------------------------------
# Per application "Known Vulnerabilities" (G) + count (H)
# ---------------------------------------------------------------------------
$knownVulnG = @{
    2  = "N/A"; 3  = "N/A"; 4  = "N/A"; 5  = "Y";   6  = "N/A"; 7  = "N/A";
    8  = "N/A"; 9  = "N/A"; 10 = "N/A"; 11 = "N/A"; 12 = "N/A"; 13 = "N/A";
    14 = "N/A"; 15 = "N/A"; 16 = "N/A"; 17 = "N/A"; 18 = "N/A"; 19 = "N/A";
    20 = "N/A"; 21 = "N/A"; 22 = "N/A"; 23 = "N/A"; 24 = "N/A"; 25 = "N/A";
    26 = "N/A"; 27 = "N/A"; 28 = "N/A"; 29 = "N/A"; 30 = "N/A"; 31 = "N/A";
}

$knownVulnH = @{
    2 = 13; 3 = 10; 4 = 14; 5 = 14; 8 = 65; 13 = 25; 21 = 65; 22 = 45; 23 = 50; 25 = 50; 29 = 25;
}

foreach ($row in 2..31) {
    $gCell = $ws.Range("G$row")
    $gCell.Value = $knownVulnG[$row]
    if ($knownVulnG[$row] -eq "Y") {
        $gCell.Style = "Good"
    } else {
        $gCell.Style = "Bad"
    }

    if ($knownVulnH.ContainsKey($row)) {
        $hCell = $ws.Range("H$row")
        $hCell.Value = $knownVulnH[$row]
        $hCell.WrapText = $true
        $hCell.HorizontalAlignment = -4108
        $hCell.VerticalAlignment = -4108
    }
}

# ---------------------------------------------------------------------------
# "Known Vulns" summary table (rows 33-44)
# ---------------------------------------------------------------------------
$ws.Range("B33").Value = "Known Vulns"
$ws.Range("B33").Style = "Bad"
$ws.Range("C33").Style = "Bad"
$ws.Range("D33").Style = "Bad"

$summary = @(
    @{ Row = 34; App = "acunetix_acuart";       Count = 13 },
    @{ Row = 35; App = "acunetix_acublog";      Count = 10 },
    @{ Row = 36; App = "acunetix_acuforum";     Count = 14 },
    @{ Row = 37; App = "altoro_mutual";         Count = 14 },
    @{ Row = 38; App = "bwapp";                 Count = 65 },
    @{ Row = 39; App = "dvwa";                  Count = 25 },
    @{ Row = 40; App = "mutillidae";            Count = 65 },
    @{ Row = 41; App = "nodegoat";              Count = 45 },
    @{ Row = 42; App = "OWASPVulnerable_app";   Count = 50 },
    @{ Row = 43; App = "security_shepherd";     Count = 50 },
    @{ Row = 44; App = "webgoat";               Count = 25 }
)

foreach ($item in $summary) {
    $r = $item.Row

    $aCell = $ws.Range("A$r")
    $aCell.Value = $item.App
    $aCell.WrapText = $true
    $aCell.VerticalAlignment = -4108

    $bCell = $ws.Range("B$r")
    $bCell.Value = $item.Count
    $bCell.WrapText = $true
    $bCell.HorizontalAlignment = -4108
    $bCell.VerticalAlignment = -4108

    $ws.Range("C$r").HorizontalAlignment = -4108
    $ws.Range("D$r").HorizontalAlignment = -4108
}

# ---------------------------------------------------------------------------
# View state: selection moved to A13, top row scrolled to row 7
# ---------------------------------------------------------------------------
$ws.Range("A13").Select()
$excel.ActiveWindow.ScrollRow = 7

Write-Host "ApplicationTracker updated"
